# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Mon Aug 26 17:28:27 UTC 2024 with GitHub Actions".
# All target cells are plain text in the source workbook (t="inlineStr"), so
# any cell whose new value looks like a bare number/decimal gets its
# NumberFormat forced to "@" (Text) first -- otherwise Excel's COM Value
# setter would auto-coerce it to a number (dropping things like trailing
# zeros, e.g. "0.590" -> 0.59).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.687.44"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.719.69"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.62"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.62"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -4.93%  "
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").Value = "3.200.49"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.66"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "63.528.33"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "2.719.75"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.19"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.64"
$ws.Range("E19").Value = "  -4.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.43"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.513"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.21"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("D28").Value = "0.0₃0879"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("E29").Value = "  +7.98%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.14"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.72"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.83"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "346.61"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("E39").Value = "  -4.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.12"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.34"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.18"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.62"
$ws.Range("E44").Value = "  -3.92%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.626"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0569"
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "132.37"
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.997"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.07"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0982"
$ws.Range("E51").Value = "  -2.95%  "
